# Update the "want to go" counts (想去人数, column F) for a handful of
# events that are listed on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8251
$ws1.Range("F5").Value = 6019
$ws1.Range("F7").Value = 100
$ws1.Range("F11").Value = 915

# Sheet "全部类型" (all types) - same events appear here on different rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8251
$ws4.Range("F5").Value = 6019
$ws4.Range("F7").Value = 100
$ws4.Range("F15").Value = 915
